$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: replace B29 / C29 / D29 contents ---

$b29 = " openActionBarOverflowOrOptionsMenu(getInstrumentation().getTargetContext());`r`n        ViewInteraction textView2 = onView(`r`n                allOf(withId(android.R.id.testo), withText(`"TestoDaSelezionare`"),`r`n                        childAtPosition(`r`n                                childAtPosition(`r`n withClassName(is(`"com.android.internal.view.menu.ListMenuItemView`")),`r`n                                        PosizioneElemento),`r`n                                Posizione),`r`n                        isDisplayed()));`r`n        textView2.perform(click());"

$c29 = "        shadowOf(activity).clickMenuItem(R.id.ElementoDelMenuDaSelezionare);"

$d29 = "solo.sendKey(solo.MENU);`r`n    // Click on Change Settings `r`n  solo.clickInList(PosizioneNelMenuDaSelezionare, 0);"

$ws.Range("B29").Value = $b29
$ws.Range("C29").Value = $c29
$ws.Range("D29").Value = $d29

# D29 gains wrap text (new style entry fontId=6 + wrapText alignment)
$ws.Range("D29").WrapText = $true

# --- Row/column sizing & selection ---
$ws.Rows.Item(29).RowHeight = 182.25
$ws.Columns.Item(2).ColumnWidth = 76.33

$ws.Range("B29").Select()
